# New PCB BOM version - replaces the rows/shared-string content of the
# BOM sheet with the refreshed component list (new sensor, decoupling caps,
# extra resistors/headers, etc.), per the "new pcb version, removed mess"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final BOM grid (header row + 22 data rows = 23 rows x up to 4 cols).
$data = @(
    @("Comment", "Designator", "Footprint", "LCSC Part Number"),
    @("100nF", "C1,C2,C3", "C_0603_1608Metric_Pad1.08x0.95mm_HandSolder", "C66501"),
    @("C_Polarized", "C4", "CP_Elec_4x5.4", "C213487"),
    @("10nF", "C5", "C_0603_1608Metric_Pad1.08x0.95mm_HandSolder", "C519406"),
    @("LEDG", "D1", "LED_0603_1608Metric", "C965806"),
    @("LEDY", "D2", "LED_0603_1608Metric", "C965802"),
    @("LEDR", "D3", "LED_0603_1608Metric", "C965798"),
    @("LEDI", "D4,D5,D6", "LED_0603_1608Metric", "C412284"),
    @("Conn_01x03_Socket", "J1,J3,J4,J7,J12,J13,J16", "PinHeader_1x03_P2.54mm_Vertical", "C49257"),
    @("Raspberry_Pi_2_3", "J2", "Raspberry_Pi_Zero_Socketed_THT_FaceDown_MountingHoles", "C2977589"),
    @("Conn_01x02_Pin", "J3_F1_BEC_OUT1", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98733"),
    @("Conn_01x02_Pin", "J3_F2_BAT1", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98733"),
    @("Conn_01x02_Pin", "J4_M0_MOT0", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98734"),
    @("Conn_01x02_Pin", "J4_M1_MOT1", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98734"),
    @("Conn_01x02_Pin", "J4_M2_MOT2", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98734"),
    @("Conn_01x02_Pin", "J4_M3_MOT3", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98734"),
    @("Conn_01x02_Pin", "J4_M4_BEC_IN1", "AMASS_XT60-F_1x02_P7.20mm_Vertical", "C98734"),
    @("Conn_01x08_Socket", "J6", "PinHeader_1x08_P2.54mm_Vertical", "C6332201"),
    @("Conn_01x03_Socket", "J14GND1", "PinHeader_1x03_P2.54mm_Vertical", "C49257"),
    @("Conn_01x03_Socket", "J15VBECOUT1", "PinHeader_1x03_P2.54mm_Vertical", "C49257"),
    @("130", "R1,R2,R3,R8,R9,R10,R11", "R_0603_1608Metric", "C22796"),
    @("49.9", "R4,R5,R6", "R_0603_1608Metric", "C23185"),
    @("LSM9DS1", "U1", "LGA-24L_3x3.5mm_P0.43mm", $null)
)

$rowCount = $data.Length
for ($r = 1; $r -le $rowCount; $r++) {
    $rowVals = $data[$r - 1]
    for ($c = 1; $c -le 4; $c++) {
        $val = $rowVals[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}

# A21/A22 ("130" / "49.9") must stay text, not be coerced to numbers - force
# text via NumberFormat, then strip the format change back off so no
# leftover style is introduced.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "130"
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "49.9"
$ws.Range("A21:A22").ClearFormats()

# Drop any leftover row-20 trailing cells beyond the new grid's extent.
$ws.Range("A24:D30").ClearContents()
